# "changes and allure correction"
# Clears out the old Test 2 - Test 5 rows (rows 7-22) leaving only the
# formatting/styles behind (matches the cells that still carry an `s`
# style attribute but no value/content in the target workbook), removes
# the now-orphaned mailto hyperlink on the old E9 cell, and resets the
# sheet selection/scroll position back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the contents (values) of the old test-case rows, keeping any
# cell formatting (fill/border/number format) that was already applied.
$ws.Range("A7:E22").ClearContents()

# Drop the hyperlink that pointed at the (now cleared) Administrator
# e-mail address in E9.
$ws.Hyperlinks.Delete()

# Scroll back to the top and select A4, instead of the old A23:E24
# selection with topLeftCell pinned at A3.
$ws.Range("A4").Select()
